# Applies the crypto price/volume update described by the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price text (dotted thousands separators, e.g. "60.799.98")
# that must stay literal text; forcing NumberFormat to Text ("@") before the
# write stops Excel from re-parsing values like "574.03" or "164.30" as
# numbers (which would also eat the trailing zero). Restoring the "Normal"
# style afterwards keeps the cell formatting identical to the original file.
function Set-TextValue($addr, $value) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

Set-TextValue 'D2' '60.799.98'
$ws.Range('E2').Value = '  -1.90%  '
Set-TextValue 'D3' '3.381.72'
$ws.Range('E3').Value = '  -0.98%  '
Set-TextValue 'D4' '0.999'
$ws.Range('E4').Value = '  -0.07%  '
Set-TextValue 'D5' '574.03'
$ws.Range('E5').Value = '  -0.71%  '
Set-TextValue 'D6' '136.52'
$ws.Range('E6').Value = '  -1.56%  '
$ws.Range('E7').Value = '  +0.02%  '
Set-TextValue 'D8' '3.379.42'
$ws.Range('E8').Value = '  -1.03%  '
$ws.Range('E9').Value = '  -1.60%  '
Set-TextValue 'D10' '7.57'
$ws.Range('E10').Value = '  +0.94%  '
$ws.Range('E11').Value = '  -3.35%  '
Set-TextValue 'D12' '0.389'
$ws.Range('E12').Value = '  -1.35%  '
Set-TextValue 'D13' '3.960.48'
$ws.Range('E13').Value = '  -1.03%  '
Set-TextValue 'D15' '26.36'
$ws.Range('E15').Value = '  +3.50%  '
$ws.Range('E16').Value = '  -3.90%  '
Set-TextValue 'D17' '3.383.74'
$ws.Range('E17').Value = '  -0.81%  '
Set-TextValue 'D18' '60.929.61'
$ws.Range('E18').Value = '  -1.77%  '
$ws.Range('E19').Value = '  -1.13%  '
Set-TextValue 'D20' '5.83'
$ws.Range('E20').Value = '  -1.05%  '
Set-TextValue 'D21' '9.46'
$ws.Range('E21').Value = '  -0.81%  '
Set-TextValue 'D22' '377.45'
$ws.Range('E22').Value = '  -3.21%  '
Set-TextValue 'D23' '0.556'
$ws.Range('E23').Value = '  -2.77%  '
Set-TextValue 'D24' '3.527.84'
$ws.Range('E24').Value = '  -0.70%  '
$ws.Range('E25').Value = '  -0.06%  '
Set-TextValue 'D26' '71.35'
$ws.Range('E26').Value = '  -0.28%  '
$ws.Range('E27').Value = '  -2.61%  '
Set-TextValue 'D28' '1.78'
$ws.Range('E28').Value = '  +11.66%  '
$ws.Range('B29').Value = 'Kaspa'
$ws.Range('C29').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
Set-TextValue 'D29' '0.169'
$ws.Range('E29').Value = '  +5.65%  '
$ws.Range('B30').Value = 'RenderToken'
$ws.Range('C30').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue 'D30' '7.56'
$ws.Range('E30').Value = '  -1.78%  '
$ws.Range('E31').Value = '  +0.16%  '
$ws.Range('E32').Value = '  -1.44%  '
$ws.Range('E33').Value = '  -0.82%  '
$ws.Range('E34').Value = '  +0.01%  '
Set-TextValue 'D35' '23.74'
$ws.Range('E35').Value = '  +0.66%  '
$ws.Range('E36').Value = '  -5.16%  '
Set-TextValue 'D37' '6.86'
$ws.Range('E37').Value = '  -2.12%  '
$ws.Range('E38').Value = '  -1.52%  '
Set-TextValue 'D39' '164.30'
$ws.Range('E39').Value = '  +0.89%  '
Set-TextValue 'D40' '0.0756'
$ws.Range('E40').Value = '  -4.64%  '
Set-TextValue 'D41' '0.999'
$ws.Range('E41').Value = '  -0.14%  '
$ws.Range('E42').Value = '  -2.31%  '
$ws.Range('E43').Value = '  -2.91%  '
$ws.Range('E44').Value = '  -1.27%  '
Set-TextValue 'D45' '41.54'
$ws.Range('E45').Value = '  -0.51%  '
$ws.Range('E46').Value = '  -2.75%  '
Set-TextValue 'D47' '24.14'
$ws.Range('E47').Value = '  -4.47%  '
Set-TextValue 'D48' '23.37'
$ws.Range('E48').Value = '  +1.20%  '
Set-TextValue 'D49' '6.81'
$ws.Range('E49').Value = '  -2.43%  '
Set-TextValue 'D50' '2.413.20'
$ws.Range('E50').Value = '  +1.52%  '
Set-TextValue 'D51' '2.40'
$ws.Range('E51').Value = '  +5.08%  '
